function Set-TextValue($cell, $val) {
    $oldStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = $oldStyle
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Simple price-only updates (column D) ---
Set-TextValue $ws.Cells.Item(2, 4) "242.38"
Set-TextValue $ws.Cells.Item(3, 4) "21.51"
Set-TextValue $ws.Cells.Item(4, 4) "5.075"
Set-TextValue $ws.Cells.Item(5, 4) "0.05593"
Set-TextValue $ws.Cells.Item(6, 4) "3.369"
Set-TextValue $ws.Cells.Item(7, 4) "6.397"
Set-TextValue $ws.Cells.Item(8, 4) "0.8050"
Set-TextValue $ws.Cells.Item(9, 4) "0.9719"
Set-TextValue $ws.Cells.Item(10, 4) "0.1407"
Set-TextValue $ws.Cells.Item(11, 4) "0.07379"
Set-TextValue $ws.Cells.Item(12, 4) "0.03100"
Set-TextValue $ws.Cells.Item(40, 4) "0.03883"
Set-TextValue $ws.Cells.Item(41, 4) "0.006928"
Set-TextValue $ws.Cells.Item(42, 4) "0.1036"
Set-TextValue $ws.Cells.Item(43, 4) "0.002913"
Set-TextValue $ws.Cells.Item(45, 4) "0.00005940"
Set-TextValue $ws.Cells.Item(49, 4) "0.09804"
Set-TextValue $ws.Cells.Item(51, 4) "0.01010"

# --- Row rotation updates (rows 13-27): Coin, Link, Price, Volume(1h) ---
$ws.Cells.Item(13, 2).Value = "BitrueCoin"
$ws.Cells.Item(13, 3).Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
Set-TextValue $ws.Cells.Item(13, 4) "0.03055"
$ws.Cells.Item(13, 5).Value = "12BitrueCoinBTR"

$ws.Cells.Item(14, 2).Value = "BitMartToken"
$ws.Cells.Item(14, 3).Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
Set-TextValue $ws.Cells.Item(14, 4) "0.09281"
$ws.Cells.Item(14, 5).Value = "13BitMartTokenBMX"

$ws.Cells.Item(15, 2).Value = "MCDex"
$ws.Cells.Item(15, 3).Value = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
Set-TextValue $ws.Cells.Item(15, 4) "3.622"
$ws.Cells.Item(15, 5).Value = "14MCDexMCB"

$ws.Cells.Item(16, 2).Value = "BitForexToken"
$ws.Cells.Item(16, 3).Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
Set-TextValue $ws.Cells.Item(16, 4) "0.001644"
$ws.Cells.Item(16, 5).Value = "15BitForexTokenBF"

$ws.Cells.Item(17, 2).Value = "CoinExToken"
$ws.Cells.Item(17, 3).Value = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
Set-TextValue $ws.Cells.Item(17, 4) "0.04704"
$ws.Cells.Item(17, 5).Value = "16CoinExTokenCET"

$ws.Cells.Item(18, 2).Value = "One"
$ws.Cells.Item(18, 3).Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
Set-TextValue $ws.Cells.Item(18, 4) "0.0005753"
$ws.Cells.Item(18, 5).Value = "17OneONE"

$ws.Cells.Item(19, 2).Value = "TigerCash"
$ws.Cells.Item(19, 3).Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
Set-TextValue $ws.Cells.Item(19, 4) "0.006394"
$ws.Cells.Item(19, 5).Value = "18TigerCashTCH"

$ws.Cells.Item(20, 2).Value = "HotbitToken"
$ws.Cells.Item(20, 3).Value = "https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb"
Set-TextValue $ws.Cells.Item(20, 4) "0.004985"
$ws.Cells.Item(20, 5).Value = "19HotbitTokenHTB"

$ws.Cells.Item(21, 2).Value = "BitKan"
$ws.Cells.Item(21, 3).Value = "https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan"
Set-TextValue $ws.Cells.Item(21, 4) "0.001043"
$ws.Cells.Item(21, 5).Value = "20BitKanKAN"

$ws.Cells.Item(22, 2).Value = "NitroEx"
$ws.Cells.Item(22, 3).Value = "https://coinranking.com/coin/8oiZw6gwYhC+nitroex-ntx"
Set-TextValue $ws.Cells.Item(22, 4) "0.0001501"
$ws.Cells.Item(22, 5).Value = "21NitroExNTX"

$ws.Cells.Item(23, 2).Value = "UpBots"
$ws.Cells.Item(23, 3).Value = "https://coinranking.com/coin/m5ozaAIK6+upbots-ubxt"
Set-TextValue $ws.Cells.Item(23, 4) "0.0003101"
$ws.Cells.Item(23, 5).Value = "22UpBotsUBXT"

$ws.Cells.Item(24, 2).Value = "LEO"
$ws.Cells.Item(24, 3).Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
Set-TextValue $ws.Cells.Item(24, 4) "3.755"
$ws.Cells.Item(24, 5).Value = "23LEOLEO"

$ws.Cells.Item(25, 2).Value = "BTSEToken"
$ws.Cells.Item(25, 3).Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
Set-TextValue $ws.Cells.Item(25, 4) "2.098"
$ws.Cells.Item(25, 5).Value = "24BTSETokenBTSE"

$ws.Cells.Item(26, 2).Value = "BitpandaEcosystemToken"
$ws.Cells.Item(26, 3).Value = "https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best"
Set-TextValue $ws.Cells.Item(26, 4) "0.3259"
$ws.Cells.Item(26, 5).Value = "25BitpandaEcosystemTokenBEST"

$ws.Cells.Item(27, 2).Value = "ProBitToken"
$ws.Cells.Item(27, 3).Value = "https://coinranking.com/coin/lQP4d6T2+probittoken-prob"
Set-TextValue $ws.Cells.Item(27, 4) "0.1287"
$ws.Cells.Item(27, 5).Value = "26ProBitTokenPROB"

